$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040482114791081
$ws.Cells.Item(2, 4).Value = 1.047968478247069
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.057082976668047
$ws.Cells.Item(2, 9).Value = 1.045060021110903
$ws.Cells.Item(2, 10).Value = 1.045568483459431
$ws.Cells.Item(2, 11).Value = 1.050729634115847
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.059818919331757
$ws.Cells.Item(2, 14).Value = 1.047053310634325
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04129269911532
$ws.Cells.Item(3, 4).Value = 1.048600439877413
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.057851590910936
$ws.Cells.Item(3, 9).Value = 1.045274921129228
$ws.Cells.Item(3, 10).Value = 1.046025397022918
$ws.Cells.Item(3, 11).Value = 1.051174013047503
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.06040141462771
$ws.Cells.Item(3, 14).Value = 1.047510873067481
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.04181762994763
$ws.Cells.Item(4, 4).Value = 1.049009711870885
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.05834962851951
$ws.Cells.Item(4, 9).Value = 1.045412879660063
$ws.Cells.Item(4, 10).Value = 1.046320794561753
$ws.Cells.Item(4, 11).Value = 1.051461215508086
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.060778356886454
$ws.Cells.Item(4, 14).Value = 1.047806690104703
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042038411774636
$ws.Cells.Item(5, 4).Value = 1.049181851938199
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.058559167404191
$ws.Cells.Item(5, 9).Value = 1.045470614160171
$ws.Cells.Item(5, 10).Value = 1.046444917104887
$ws.Cells.Item(5, 11).Value = 1.051581872596618
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.060936828698847
$ws.Cells.Item(5, 14).Value = 1.047930988916086
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042075487882716
$ws.Cells.Item(6, 4).Value = 1.049210759776076
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.058594359456888
$ws.Cells.Item(6, 9).Value = 1.045480292574153
$ws.Cells.Item(6, 10).Value = 1.046465754119141
$ws.Cells.Item(6, 11).Value = 1.051602126557112
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.060963437076331
$ws.Cells.Item(6, 14).Value = 1.04795185552129
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041820579650844
$ws.Cells.Item(7, 4).Value = 1.049012011694481
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.058352427746925
$ws.Cells.Item(7, 9).Value = 1.045413652146793
$ws.Cells.Item(7, 10).Value = 1.046322453340754
$ws.Cells.Item(7, 11).Value = 1.051462828061127
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.060780474375832
$ws.Cells.Item(7, 14).Value = 1.047808351239361
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04075596579996
$ws.Cells.Item(8, 4).Value = 1.048181978969573
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.057342588810166
$ws.Cells.Item(8, 9).Value = 1.045132874035101
$ws.Cells.Item(8, 10).Value = 1.045722951796999
$ws.Cells.Item(8, 11).Value = 1.050879883745224
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.060015769396599
$ws.Cells.Item(8, 14).Value = 1.04720799833465
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038883337504817
$ws.Cells.Item(9, 4).Value = 1.046722119766167
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.055568517184462
$ws.Cells.Item(9, 9).Value = 1.044629749300436
$ws.Cells.Item(9, 10).Value = 1.044664648403484
$ws.Cells.Item(9, 11).Value = 1.049850112795776
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.05866855406329
$ws.Cells.Item(9, 14).Value = 1.046148192028927
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037637276722903
$ws.Cells.Item(10, 4).Value = 1.045750846388047
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.05438954724659
$ws.Cells.Item(10, 9).Value = 1.044288763789167
$ws.Cells.Item(10, 10).Value = 1.043957907686101
$ws.Cells.Item(10, 11).Value = 1.049161964304343
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.057770703649912
$ws.Cells.Item(10, 14).Value = 1.045440447658661
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037098299119835
$ws.Cells.Item(11, 4).Value = 1.045330763291816
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.053879953318597
$ws.Cells.Item(11, 9).Value = 1.04413980296866
$ws.Cells.Item(11, 10).Value = 1.043651611251195
$ws.Cells.Item(11, 11).Value = 1.048863617358357
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.057382014369193
$ws.Cells.Item(11, 14).Value = 1.045133716247685
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036898186912225
$ws.Cells.Item(12, 4).Value = 1.045174800377008
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.053690805853054
$ws.Cells.Item(12, 9).Value = 1.044084276030162
$ws.Cells.Item(12, 10).Value = 1.04353779937151
$ws.Cells.Item(12, 11).Value = 1.048752743145762
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.057237652401309
$ws.Cells.Item(12, 14).Value = 1.045019742742076
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036941107647392
$ws.Cells.Item(13, 4).Value = 1.045208251544147
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.053731372330837
$ws.Cells.Item(13, 9).Value = 1.044096195614643
$ws.Cells.Item(13, 10).Value = 1.043562214184167
$ws.Cells.Item(13, 11).Value = 1.048776528502472
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.057268617849378
$ws.Cells.Item(13, 14).Value = 1.045044192226567
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037081755975882
$ws.Cells.Item(14, 4).Value = 1.045317869815176
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.053864315494936
$ws.Cells.Item(14, 9).Value = 1.044135217093374
$ws.Cells.Item(14, 10).Value = 1.04364220432915
$ws.Cells.Item(14, 11).Value = 1.048854453573048
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.057370081057883
$ws.Cells.Item(14, 14).Value = 1.045124295966732
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037168425774158
$ws.Cells.Item(15, 4).Value = 1.045385419196598
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.053946244566831
$ws.Cells.Item(15, 9).Value = 1.044159233533211
$ws.Cells.Item(15, 10).Value = 1.043691483668894
$ws.Cells.Item(15, 11).Value = 1.048902458549517
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.057432597864592
$ws.Cells.Item(15, 14).Value = 1.04517364528879
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037673059149211
$ws.Cells.Item(16, 4).Value = 1.045778736287936
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.054423386616562
$ws.Cells.Item(16, 9).Value = 1.044298622260227
$ws.Cells.Item(16, 10).Value = 1.043978229927956
$ws.Cells.Item(16, 11).Value = 1.049181756840315
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.057796501626287
$ws.Cells.Item(16, 14).Value = 1.04546079876043
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037989757456034
$ws.Cells.Item(17, 4).Value = 1.046025584767176
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.054722929696439
$ws.Cells.Item(17, 9).Value = 1.044385706495606
$ws.Cells.Item(17, 10).Value = 1.044158026080428
$ws.Cells.Item(17, 11).Value = 1.049356854060931
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.058024792987905
$ws.Cells.Item(17, 14).Value = 1.045640850244062
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038174537579826
$ws.Cells.Item(18, 4).Value = 1.046169613951233
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.054897735564883
$ws.Cells.Item(18, 9).Value = 1.044436374626968
$ws.Cells.Item(18, 10).Value = 1.044262871704335
$ws.Cells.Item(18, 11).Value = 1.049458949122012
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.058157959652671
$ws.Cells.Item(18, 14).Value = 1.045745844760781
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038237552165304
$ws.Cells.Item(19, 4).Value = 1.046218732036507
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.054957354621591
$ws.Cells.Item(19, 9).Value = 1.044453629646458
$ws.Cells.Item(19, 10).Value = 1.044298616830847
$ws.Cells.Item(19, 11).Value = 1.049493754721254
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.058203367381002
$ws.Cells.Item(19, 14).Value = 1.045781640649473
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037955772986587
$ws.Cells.Item(20, 4).Value = 1.045999095423828
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.054690782498086
$ws.Cells.Item(20, 9).Value = 1.044376376273212
$ws.Cells.Item(20, 10).Value = 1.044138738379405
$ws.Cells.Item(20, 11).Value = 1.049338071519457
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.05800029862785
$ws.Cells.Item(20, 14).Value = 1.045621535152291
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037040336103874
$ws.Cells.Item(21, 4).Value = 1.045285587891079
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.05382516320968
$ws.Cells.Item(21, 9).Value = 1.044123731651445
$ws.Cells.Item(21, 10).Value = 1.043618650312098
$ws.Cells.Item(21, 11).Value = 1.048831508086062
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.057340202257451
$ws.Cells.Item(21, 14).Value = 1.045100708500273
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036465275000007
$ws.Cells.Item(22, 4).Value = 1.04483741020035
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.053281715794819
$ws.Cells.Item(22, 9).Value = 1.043963748878868
$ws.Cells.Item(22, 10).Value = 1.043291420975026
$ws.Cells.Item(22, 11).Value = 1.048512694730687
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.05692525838737
$ws.Cells.Item(22, 14).Value = 1.044773014460008
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036770076825805
$ws.Cells.Item(23, 4).Value = 1.04507495603721
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.053569730878777
$ws.Cells.Item(23, 9).Value = 1.044048666087062
$ws.Cells.Item(23, 10).Value = 1.04346491274991
$ws.Cells.Item(23, 11).Value = 1.048681733308722
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.057145219388977
$ws.Cells.Item(23, 14).Value = 1.044946752613114
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037971128934636
$ws.Cells.Item(24, 4).Value = 1.046011064673984
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.05470530816329
$ws.Cells.Item(24, 9).Value = 1.04438059259046
$ws.Cells.Item(24, 10).Value = 1.044147453743088
$ws.Cells.Item(24, 11).Value = 1.049346558653435
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.058011366549753
$ws.Cells.Item(24, 14).Value = 1.04563026289279
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.039367048561788
$ws.Cells.Item(25, 4).Value = 1.047099189675551
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.056026505872452
$ws.Cells.Item(25, 9).Value = 1.044760804368306
$ws.Cells.Item(25, 10).Value = 1.044938463151115
$ws.Cells.Item(25, 11).Value = 1.050116627144756
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.059016796906102
$ws.Cells.Item(25, 14).Value = 1.047053310634325
